$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 31   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/5/2024  Through  2/11/2024"

# --- Cells whose type/style changes (text <-> number) ---
# Copy number-format from a same-column reference cell, then set the new value.
# For text values that look numeric ("0"), prefix with an apostrophe to force text,
# then re-apply the correct (non quote-prefixed) style via PasteSpecial formats.
$ws.Range("D20").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("C16").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1

$ws.Range("D16").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1

$ws.Range("E16").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = -100

$ws.Range("G16").Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("G26").Value = 1

$ws.Range("H16").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("H26").Value = -100

$ws.Range("C16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 3

$ws.Range("C30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Plain value updates (style/type unchanged) ---
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -11.111111111111
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 21
$ws.Range("K16").Value = -38.095238095238
$ws.Range("L16").Value = -40.909090909090
$ws.Range("M16").Value = -40.909090909090
$ws.Range("N16").Value = -87.378640776699
$ws.Range("D17").Value = 4
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -77.777777777777
$ws.Range("J17").Value = 16
$ws.Range("K17").Value = -50
$ws.Range("N17").Value = -65.217391304347
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -30
$ws.Range("I18").Value = 25
$ws.Range("J18").Value = 38
$ws.Range("K18").Value = -34.210526315789
$ws.Range("L18").Value = -21.875
$ws.Range("M18").Value = -10.714285714285
$ws.Range("N18").Value = -68.75
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = -8.695652173913
$ws.Range("F19").Value = 76
$ws.Range("G19").Value = 101
$ws.Range("H19").Value = -24.752475247524
$ws.Range("I19").Value = 115
$ws.Range("J19").Value = 155
$ws.Range("K19").Value = -25.806451612903
$ws.Range("L19").Value = 6.481481481481
$ws.Range("M19").Value = 4.545454545454
$ws.Range("N19").Value = -56.439393939393
$ws.Range("L20").Value = -57.142857142857
$ws.Range("M20").Value = -40
$ws.Range("N20").Value = -97
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -29.411764705882
$ws.Range("F21").Value = 101
$ws.Range("G21").Value = 141
$ws.Range("H21").Value = -28.368794326241
$ws.Range("I21").Value = 164
$ws.Range("J21").Value = 234
$ws.Range("K21").Value = -29.914529914529
$ws.Range("L21").Value = -8.379888268156
$ws.Range("M21").Value = -5.202312138728
$ws.Range("N21").Value = -71.278458844133
$ws.Range("E22").Value = -50
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 6
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = -14.285714285714
$ws.Range("M22").Value = -53.846153846153
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 129
$ws.Range("G24").Value = 152
$ws.Range("H24").Value = -15.131578947368
$ws.Range("I24").Value = 180
$ws.Range("J24").Value = 213
$ws.Range("K24").Value = -15.492957746478
$ws.Range("L24").Value = -5.263157894736
$ws.Range("M24").Value = 7.142857142857
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -42.857142857142
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = -17.857142857142
$ws.Range("I25").Value = 40
$ws.Range("J25").Value = 44
$ws.Range("K25").Value = -9.090909090909
$ws.Range("L25").Value = 122.222222222222
$ws.Range("M25").Value = 135.294117647059
$ws.Range("J26").Value = 2
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = -50
$ws.Range("L27").Value = 25
